# Updated cryptos list on Sat Jan  6 23:43:42 UTC 2024 with GitHub Actions
# Applies the refreshed Price/Volume(1h) figures (and a handful of reordered
# rows) from the upstream coinranking.com snapshot to the "cryptos" sheet.
# Numeric-looking "Price" strings (e.g. "1.00", "0.523") are prefixed with a
# leading apostrophe so Excel keeps storing/display them as text, exactly as
# they were before (matching the source inlineStr cell type) instead of
# silently re-parsing them as numbers and dropping formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.306.92"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "2.239.76"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'307.26"
$ws.Range("E5").Value = "  -3.35%  "

$ws.Range("D6").Value = "'93.51"
$ws.Range("E6").Value = "  -6.57%  "

$ws.Range("E7").Value = "  -1.24%  "

$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  -3.10%  "

$ws.Range("D10").Value = "'34.39"
$ws.Range("E10").Value = "  -5.51%  "

$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = "  -2.61%  "

$ws.Range("D12").Value = "'7.14"
$ws.Range("E12").Value = "  -4.28%  "

$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").Value = "2.329.50"
$ws.Range("E14").Value = "  +2.92%  "

$ws.Range("D15").Value = "'0.828"
$ws.Range("E15").Value = "  -2.66%  "

$ws.Range("D16").Value = "'13.42"
$ws.Range("E16").Value = "  -4.63%  "

$ws.Range("D17").Value = "44.048.85"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("E18").Value = "  -1.98%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = "  -1.53%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.01"
$ws.Range("E20").Value = "  -9.23%  "

$ws.Range("D21").Value = "'65.62"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'236.90"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "'2.94"
$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("E24").Value = "  -2.69%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").Value = "'39.64"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("D27").Value = "'2.21"
$ws.Range("E27").Value = "  +3.87%  "

$ws.Range("D28").Value = "'9.84"
$ws.Range("E28").Value = "  -4.12%  "

$ws.Range("D29").Value = "'20.01"
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").Value = "'5.88"
$ws.Range("E30").Value = "  -3.72%  "

$ws.Range("D31").Value = "'151.78"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("E32").Value = "  -6.28%  "

$ws.Range("E33").Value = "  -2.88%  "

$ws.Range("E34").Value = "  -12.72%  "

$ws.Range("E35").Value = "  +0.84%  "

$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("E37").Value = "  -8.92%  "

$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'3.78"
$ws.Range("E39").Value = "  -5.80%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'14.23"
$ws.Range("E40").Value = "  -8.63%  "

$ws.Range("E41").Value = "  -4.26%  "

$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").Value = "1.700.86"
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("D44").Value = "'82.32"
$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("D45").Value = "'0.191"
$ws.Range("E45").Value = "  -3.37%  "

$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "'4.92"
$ws.Range("E46").Value = "  -5.88%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'99.22"
$ws.Range("E47").Value = "  -3.49%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.07"
$ws.Range("E49").Value = "  -2.65%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'54.70"
$ws.Range("E50").Value = "  -4.52%  "

$ws.Range("D51").Value = "'67.28"
$ws.Range("E51").Value = "  -5.99%  "
